# Update column C (predicted factor values) with re-computed results.
# The underlying loop that generates these numbers was fixed/adjusted,
# so the outputs for rows 2-11 in column C now differ slightly from
# the previous run, while columns A and B remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = -5.023377998469892
$ws.Range("C3").Value  = -1.151483576865103
$ws.Range("C4").Value  = -0.06510371792358649
$ws.Range("C5").Value  = -0.4130396773035255
$ws.Range("C6").Value  = 0.01420778898168975
$ws.Range("C7").Value  = 0.1031345171695419
$ws.Range("C8").Value  = 0.1337589614876022
$ws.Range("C9").Value  = 0.02733475124857185
$ws.Range("C10").Value = 0.0234479315146108
$ws.Range("C11").Value = 0.006722294386397641
